{"js": "// 1. Fix the duplicated wording: \"the player must player and kill enemies\" -> \"the player must kill enemies\".\nconst dupFix = context.document.body.search(\"the player must player and kill enemies\", { matchCase: true });\ndupFix.load(\"items\");\nawait context.sync();\nif (dupFix.items.length > 0) {\n  dupFix.items[0].insertText(\"the player must kill enemies\", \"Replace\");\n}\nawait context.sync();\n\n// 2. The paragraph currently reads (in document order):\n//      \"...use of a \\u201cdash\\u201d. \" [bookmark] \"Once all enemies ... the next\" \", with the difficulty ... is hit.\"\n//    The two sentences that trail the _GoBack bookmark need to move in front of it, and a brand new\n//    closing sentence about the player's score needs to be appended after them - with the bookmark\n//    ending up at the very end of the paragraph. Grab that trailing text via the bookmark itself so\n//    this does not depend on matching the (fairly generic) sentence text more than once.\nconst body = context.document.body;\nconst bookmark = body.getBookmarkRangeOrNullObject(\"_GoBack\");\nbookmark.load(\"isNullObject\");\nawait context.sync();\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nconst paragraph = paragraphs.items[1];\n\nconst trailingRange = bookmark.getRange(\"After\").expandTo(paragraph.getRange(\"End\"));\ntrailingRange.load(\"text\");\nawait context.sync();\nconst trailingText = trailingRange.text; // \"Once all enemies ... the player is hit.\"\n\n// 3. Remove that trailing text; it gets reinstated (in front of the bookmark) in the next step.\ntrailingRange.insertText(\"\", \"Replace\");\nawait context.sync();\n\n// 4. Locate the \". \" that used to end the paragraph's first sentence (right after the dash-quote)\n//    and replace it with \". \" followed by the two sentences we just pulled out, then tack on the\n//    brand new closing sentence about the player's score.\nconst sentenceGap = paragraph.search(\". \", { matchCase: true });\nsentenceGap.load(\"items\");\nawait context.sync();\nlet cursor = sentenceGap.items[0].insertText(\". \" + trailingText, \"Replace\");\nawait context.sync();\n\ncursor = cursor.getRange(\"After\").insertText(\" The player\\u2019s score is determined by the number of enemies killed.\", \"Replace\");\nawait context.sync();\n\n// 5. Move the _GoBack bookmark so that it sits at the very end of the paragraph, after all the text.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst paragraphsAfter = body.paragraphs;\nparagraphsAfter.load(\"items\");\nawait context.sync();\nparagraphsAfter.items[1].getRange(\"End\").insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Fix the duplicated wording: \"the player must player and kill enemies\" -> \"the player must kill enemies\".\n$fixRange = $d.Content\n$fixRange.Find.Execute(\"the player must player and kill enemies\", $false, $false, $false, $false, $false, $true, 1, $false, \"the player must kill enemies\", 2) | Out-Null\n\n# 2. The paragraph currently reads (in document order):\n#      \"...use of a \"dash\". \" [bookmark] \"Once all enemies ... the next\" \", with the difficulty ... is hit.\"\n#    The two sentences trailing the _GoBack bookmark need to move in front of it, and a brand new\n#    closing sentence about the player's score needs to be appended after them, with the bookmark\n#    ending up at the very end of the paragraph.\n$para = $d.Paragraphs(2).Range\n\n$bookmark = $d.Bookmarks(\"_GoBack\")\n$bookmarkRange = $bookmark.Range\n\n# Text currently trailing the bookmark, up to (but excluding) the paragraph mark.\n$leftoverRange = $d.Range($bookmarkRange.End, $para.End - 1)\n$leftoverText = $leftoverRange.Text\n\n# 3. Remove that trailing text; it gets reinstated (in front of the bookmark) below.\n$leftoverRange.Text = \"\"\n\n# 4. Locate the \". \" that ends the paragraph's first sentence (right after the dash-quote) and\n#    replace it with \". \" followed by the two sentences pulled out above, then tack on the brand\n#    new closing sentence about the player's score.\n$gapRange = $d.Paragraphs(2).Range.Duplicate\n$gapRange.Find.Execute(\". \") | Out-Null\n$gapRange.Text = \". \" + $leftoverText + \" The player\" + [char]0x2019 + \"s score is determined by the number of enemies killed.\"\n\n# 5. Move the _GoBack bookmark so that it sits at the very end of the paragraph, after all the text.\n#    Adding a bookmark with a zero-length range that lands exactly on a paragraph-end boundary is\n#    placed incorrectly, so park a throw-away sentinel character there first, anchor the (now\n#    interior, not boundary) collapsed bookmark immediately before it, then remove the sentinel.\n$d.Bookmarks(\"_GoBack\").Delete()\n$endOfPara = $d.Paragraphs(2).Range.Duplicate\n$endOfPara.MoveEnd(1, -1)\n$endOfPara.Collapse(0)\n$endOfPara.InsertAfter(\"X\")\n\n$paraWithSentinel = $d.Paragraphs(2).Range.Duplicate\n$paraWithSentinel.MoveEnd(1, -1)\n$bookmarkSpot = $d.Range($paraWithSentinel.End - 1, $paraWithSentinel.End - 1)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkSpot) | Out-Null\n\n$sentinelRange = $d.Range($d.Bookmarks(\"_GoBack\").Range.End, $d.Bookmarks(\"_GoBack\").Range.End + 1)\n$sentinelRange.Text = \"\"\n"}
